$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data to reflect latest scrape
# Leading apostrophe forces Excel to store the value as literal text,
# preserving the General number format and avoiding numeric/date coercion
# of values such as "1.001", "22.04", "0.00001116", etc.

$ws.Range("D2").Value = "'27.430.52"
$ws.Range("E2").Value = "'  +9.73%  "

$ws.Range("D3").Value = "'1.789.35"
$ws.Range("E3").Value = "'  +7.11%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "'  -0.18%  "

$ws.Range("D5").Value = "'341.45"
$ws.Range("E5").Value = "'  +3.90%  "

$ws.Range("D6").Value = "'0.9987"
$ws.Range("E6").Value = "'  -0.19%  "

$ws.Range("D7").Value = "'0.3792"
$ws.Range("E7").Value = "'  +4.17%  "

$ws.Range("E8").Value = "'  +8.33%  "

$ws.Range("D9").Value = "'49.44"
$ws.Range("E9").Value = "'  +4.79%  "

$ws.Range("D10").Value = "'1.218"
$ws.Range("E10").Value = "'  +6.55%  "

$ws.Range("D11").Value = "'0.07712"
$ws.Range("E11").Value = "'  +6.99%  "

$ws.Range("D12").Value = "'0.9998"
$ws.Range("E12").Value = "'  -0.10%  "

$ws.Range("B13").Value = "'Solana"
$ws.Range("C13").Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "'22.04"
$ws.Range("E13").Value = "'  +12.65%  "

$ws.Range("B14").Value = "'Polkadot"
$ws.Range("C14").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.612"
$ws.Range("E14").Value = "'  +9.05%  "

$ws.Range("D15").Value = "'7.189"
$ws.Range("E15").Value = "'  +8.19%  "

$ws.Range("D16").Value = "'1.784.46"
$ws.Range("E16").Value = "'  +6.66%  "

$ws.Range("D17").Value = "'0.00001116"
$ws.Range("E17").Value = "'  +6.28%  "

$ws.Range("D18").Value = "'0.06788"
$ws.Range("E18").Value = "'  +4.02%  "

$ws.Range("D19").Value = "'85.77"
$ws.Range("E19").Value = "'  +8.90%  "

$ws.Range("D20").Value = "'0.9986"
$ws.Range("E20").Value = "'  -0.18%  "

$ws.Range("D21").Value = "'17.75"
$ws.Range("E21").Value = "'  +12.60%  "

$ws.Range("D22").Value = "'6.426"
$ws.Range("E22").Value = "'  +9.30%  "

$ws.Range("D23").Value = "'13.07"
$ws.Range("E23").Value = "'  +2.36%  "

$ws.Range("D24").Value = "'27.419.98"
$ws.Range("E24").Value = "'  +9.68%  "

$ws.Range("D25").Value = "'2.463"
$ws.Range("E25").Value = "'  +1.05%  "

$ws.Range("B26").Value = "'LidoDAOToken"
$ws.Range("C26").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.555"
$ws.Range("E26").Value = "'  +7.74%  "

$ws.Range("B27").Value = "'ImmutableX"
$ws.Range("C27").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "'1.494"
$ws.Range("E27").Value = "'  +26.30%  "

$ws.Range("D28").Value = "'20.72"
$ws.Range("E28").Value = "'  +11.18%  "

$ws.Range("D29").Value = "'153.93"
$ws.Range("E29").Value = "'  +2.89%  "

$ws.Range("D30").Value = "'1.980.64"
$ws.Range("E30").Value = "'  +6.50%  "

$ws.Range("D31").Value = "'136.35"
$ws.Range("E31").Value = "'  +8.17%  "

$ws.Range("B32").Value = "'HuobiToken"
$ws.Range("C32").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").Value = "'4.187"
$ws.Range("E32").Value = "'  +2.75%  "

$ws.Range("B33").Value = "'Filecoin"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'6.303"
$ws.Range("E33").Value = "'  +9.65%  "

$ws.Range("B34").Value = "'Aptos"
$ws.Range("C34").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "'13.85"
$ws.Range("E34").Value = "'  +13.47%  "

$ws.Range("B35").Value = "'Stellar"
$ws.Range("C35").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "'0.08793"
$ws.Range("E35").Value = "'  +4.58%  "

$ws.Range("D36").Value = "'1.716"
$ws.Range("E36").Value = "'  +4.02%  "

$ws.Range("D37").Value = "'5.641"
$ws.Range("E37").Value = "'  +9.55%  "

$ws.Range("B38").Value = "'Hedera"
$ws.Range("C38").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06537"
$ws.Range("E38").Value = "'  +8.02%  "

$ws.Range("B39").Value = "'VeChain"
$ws.Range("C39").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02422"
$ws.Range("E39").Value = "'  +9.13%  "

$ws.Range("B40").Value = "'Algorand"
$ws.Range("C40").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2255"
$ws.Range("E40").Value = "'  +8.60%  "

$ws.Range("D41").Value = "'0.6814"
$ws.Range("E41").Value = "'  +14.94%  "

$ws.Range("D42").Value = "'8.947"
$ws.Range("E42").Value = "'  +8.06%  "

$ws.Range("D43").Value = "'1.246"
$ws.Range("E43").Value = "'  +1.02%  "

$ws.Range("D44").Value = "'14.70"
$ws.Range("E44").Value = "'  +8.91%  "

$ws.Range("D45").Value = "'0.6443"
$ws.Range("E45").Value = "'  +13.38%  "

$ws.Range("D46").Value = "'0.9984"
$ws.Range("E46").Value = "'  -0.19%  "

$ws.Range("D47").Value = "'4.014"
$ws.Range("E47").Value = "'  +4.93%  "

$ws.Range("D48").Value = "'2.180"
$ws.Range("E48").Value = "'  +11.77%  "

$ws.Range("D49").Value = "'132.65"
$ws.Range("E49").Value = "'  +6.90%  "

$ws.Range("D50").Value = "'0.07349"
$ws.Range("E50").Value = "'  +3.04%  "

$ws.Range("D51").Value = "'80.67"
$ws.Range("E51").Value = "'  +8.65%  "
